# Insert a new weekly record (2 rows: "Primera" and "Segunda" quality)
# right before the existing row 430 of the "Betarraga" data table, pushing
# every row from the old 430 onward down by two rows (old 430 -> new 432,
# ..., old 512 -> new 514). This matches the diff: dimension grows from
# A1:R512 to A1:R514, and all the date/price values seen shifting by two
# rows are simply the pre-existing rows sliding down to make room for the
# two brand-new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 430:431 - everything currently at row 430 and
# below (through 512) shifts down to 432..514.
$ws.Rows("430:431").Insert()

# New row 430 - "Primera" quality record for date 44694 (2022-05-13)
$ws.Cells.Item(430, 1).Value  = 9
$ws.Cells.Item(430, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(430, 3).Value  = "Metropolitana"
$ws.Cells.Item(430, 4).Value  = 44694
$ws.Cells.Item(430, 5).Value  = 13
$ws.Cells.Item(430, 6).Value  = 100114014
$ws.Cells.Item(430, 7).Value  = "Betarraga"
$ws.Cells.Item(430, 8).Value  = "Sin especificar"
$ws.Cells.Item(430, 9).Value  = "Primera"
$ws.Cells.Item(430, 10).Value = 9700
$ws.Cells.Item(430, 11).Value = 110
$ws.Cells.Item(430, 12).Value = 120
$ws.Cells.Item(430, 13).Value = 115
$ws.Cells.Item(430, 14).Value = "$/unidad"
$ws.Cells.Item(430, 15).Value = "Región Metropolitana"
$ws.Cells.Item(430, 16).Value = 115
$ws.Cells.Item(430, 17).Value = 1
$ws.Cells.Item(430, 18).Value = "Hortaliza"

# New row 431 - "Segunda" quality record for the same date 44694
$ws.Cells.Item(431, 1).Value  = 9
$ws.Cells.Item(431, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(431, 3).Value  = "Metropolitana"
$ws.Cells.Item(431, 4).Value  = 44694
$ws.Cells.Item(431, 5).Value  = 13
$ws.Cells.Item(431, 6).Value  = 100114014
$ws.Cells.Item(431, 7).Value  = "Betarraga"
$ws.Cells.Item(431, 8).Value  = "Sin especificar"
$ws.Cells.Item(431, 9).Value  = "Segunda"
$ws.Cells.Item(431, 10).Value = 3400
$ws.Cells.Item(431, 11).Value = 100
$ws.Cells.Item(431, 12).Value = 100
$ws.Cells.Item(431, 13).Value = 100
$ws.Cells.Item(431, 14).Value = "$/unidad"
$ws.Cells.Item(431, 15).Value = "Región Metropolitana"
$ws.Cells.Item(431, 16).Value = 100
$ws.Cells.Item(431, 17).Value = 1
$ws.Cells.Item(431, 18).Value = "Hortaliza"
